$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) values - use leading apostrophe to force text
# entry (matching the source data which stores these as plain text), then
# clear the auto-applied "Text" number format so the cell keeps the default style.
$ws.Range('D2').Value = '''65.694.49'
$ws.Range('D2').ClearFormats()
$ws.Range('D3').Value = '''3.271.95'
$ws.Range('D3').ClearFormats()
$ws.Range('D5').Value = '''582.52'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').Value = '''178.78'
$ws.Range('D6').ClearFormats()
$ws.Range('D10').Value = '''6.73'
$ws.Range('D10').ClearFormats()
$ws.Range('D12').Value = '''3.845.94'
$ws.Range('D12').ClearFormats()
$ws.Range('D14').Value = '''65.840.70'
$ws.Range('D14').ClearFormats()
$ws.Range('D15').Value = '''25.93'
$ws.Range('D15').ClearFormats()
$ws.Range('D17').Value = '''3.220.95'
$ws.Range('D17').ClearFormats()
$ws.Range('D18').Value = '''425.11'
$ws.Range('D18').ClearFormats()
$ws.Range('D19').Value = '''13.14'
$ws.Range('D19').ClearFormats()
$ws.Range('D22').Value = '''71.73'
$ws.Range('D22').ClearFormats()
$ws.Range('D25').Value = '''3.422.37'
$ws.Range('D25').ClearFormats()
$ws.Range('D31').Value = '''1.94'
$ws.Range('D31').ClearFormats()
$ws.Range('D32').Value = '''22.15'
$ws.Range('D32').ClearFormats()
$ws.Range('D33').Value = '''1.00'
$ws.Range('D33').ClearFormats()
$ws.Range('D37').Value = '''159.25'
$ws.Range('D37').ClearFormats()
$ws.Range('D40').Value = '''26.30'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').Value = '''2.780.13'
$ws.Range('D41').ClearFormats()
$ws.Range('D42').Value = '''0.764'
$ws.Range('D42').ClearFormats()
$ws.Range('D43').Value = '''4.31'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').Value = '''39.93'
$ws.Range('D44').ClearFormats()
$ws.Range('D46').Value = '''5.84'
$ws.Range('D46').ClearFormats()
$ws.Range('D48').Value = '''314.11'
$ws.Range('D48').ClearFormats()

# Update "Volume(1h)" column (E) values (plain text assignment is sufficient
# since these contain "%" and spaces, which Excel always treats as text).
$ws.Range('E2').Value = '  -1.60%  '
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +1.59%  '
$ws.Range('E6').Value = '  -2.64%  '
$ws.Range('E7').Value = '  +7.15%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -4.10%  '
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('E13').Value = '  -4.73%  '
$ws.Range('E14').Value = '  -1.43%  '
$ws.Range('E16').Value = '  -3.42%  '
$ws.Range('E17').Value = '  -3.23%  '
$ws.Range('E18').Value = '  -2.11%  '
$ws.Range('E19').Value = '  -4.48%  '
$ws.Range('E20').Value = '  -3.69%  '
$ws.Range('E21').Value = '  -3.90%  '
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('E25').Value = '  -1.05%  '
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('E28').Value = '  -5.67%  '
$ws.Range('E29').Value = '  -2.54%  '
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('E32').Value = '  -2.94%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  -4.19%  '
$ws.Range('E35').Value = '  -3.45%  '
$ws.Range('E36').Value = '  -4.54%  '
$ws.Range('E37').Value = '  -0.69%  '
$ws.Range('E38').Value = '  -6.24%  '
$ws.Range('E39').Value = '  -3.42%  '
$ws.Range('E40').Value = '  -3.89%  '
$ws.Range('E41').Value = '  -1.57%  '
$ws.Range('E42').Value = '  -3.51%  '
$ws.Range('E43').Value = '  -3.34%  '
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('E45').Value = '  -3.26%  '
$ws.Range('E46').Value = '  -6.88%  '
$ws.Range('E47').Value = '  -3.69%  '
$ws.Range('E48').Value = '  -2.27%  '
$ws.Range('E49').Value = '  -6.20%  '
$ws.Range('E50').Value = '  -2.79%  '
$ws.Range('E51').Value = '  +4.29%  '
